$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.267.17"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "1.867.03"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4370"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3702"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9404"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("D12").Value = "1.931.99"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.726"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.455"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06864"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "82.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009109"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").Value = "28.260.15"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.132"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "2.117.79"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.019"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.326"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.730"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09034"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8012"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.846"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.171"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.954"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05442"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01954"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("E40").Value = "  +7.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.138"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5248"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1673"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.704"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06777"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.047"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4879"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000002530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.679"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.99%  "
